# manage.xlsx — "starting with ptpiree conv (model and controller)"
#
# Adds a new worksheet "ptpiree conv" after the last existing sheet
# (tasks, tests, log and auth, invoices) and fills it in with the
# model/controller planning notes, becoming the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- add the new sheet after the current last sheet (-> becomes sheet5,
#     last tab, and is made the active sheet/tab automatically) ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ptpiree conv"

# --- fill in the cells --------------------------------------------------
# Order matters here: it reproduces the shared-string insertion order of
# the authored workbook (models: / ptpiree_content / head / data / date /
# values / obj to file / methods: / How does it work / ...).
$ws.Range("J1").Value = "models:"
$ws.Range("K2").Value = "ptpiree_content"
$ws.Range("L2").Value = "head"
$ws.Range("M2").Value = "data"
$ws.Range("N2").Value = "date"
$ws.Range("O2").Value = "values"
$ws.Range("H2").Value = "obj to file"
$ws.Range("G1").Value = "methods:"
$ws.Range("D1").Value = "How does it work"
$ws.Range("D3").Value = "1. import csv file with full content (1line = 1output file)"
$ws.Range("D2").Value = "full auto or high managable"
$ws.Range("D4").Value = "2. convert every line into obj"
$ws.Range("D5").Value = "3. add every obj to array"
$ws.Range("D6").Value = "4. conv obj to file"
$ws.Range("D7").Value = "5. add to arch or other array"
$ws.Range("D8").Value = "6. download"
$ws.Range("E3").Value = "send via post/get/other"
$ws.Range("D10").Value = "OR easy way line by line"
$ws.Range("D11").Value = "1. conv string to array"
$ws.Range("D12").Value = "2. array pos to file line"

# --- column widths (best-fit-like widths from the authored sheet) ------
$ws.Columns.Item(4).ColumnWidth = 50.7109375
$ws.Columns.Item(7).ColumnWidth = 9.42578125
$ws.Columns.Item(8).ColumnWidth = 9.5703125
$ws.Columns.Item(9).ColumnWidth = 2.7109375
$ws.Columns.Item(11).ColumnWidth = 15.5703125

# --- selection on the new sheet, as left by the author ------------------
$ws.Range("D19").Select() | Out-Null
